$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: title/link update
$ws.Range("D9").Value = "[공지] SIAI 지원자 파비 장학금"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-siai-pabii-grant/#utm_source=rss&utm_medium=rss&utm_campaign=notice-siai-pabii-grant"

# Row 42: title/link update
$ws.Range("D42").Value = "Python 네이버 밴드 - 키워드 분석"
$ws.Range("E42").Value = "https://kjk92.tistory.com/91"

# Row 51: title/link update
$ws.Range("D51").Value = "[MySQL/MariaDB] 테이블 생성 쿼리 확인"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MYSQL-MariaDB-%ED%85%8C%EC%9D%B4%EB%B8%94-%EC%83%9D%EC%84%B1-%EC%BF%BC%EB%A6%AC-%ED%99%95%EC%9D%B8"

# Row 52: title update only
$ws.Range("D52").Value = "숨은 DS"
